$d = $word.ActiveDocument
$word.Options.StoreRSIDOnSave = $false

# Find the paragraph that ends with the italics sentence, then insert two new
# paragraphs after it: a FirstParagraph-styled one, and a BlockText-styled one.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*I don*t know why this would be in italics*") {
        $target = $p
    }
}

$r = $target.Range
$r.InsertParagraphAfter()

# Re-enumerate to find the newly inserted (empty) paragraph right after target.
$newPara1 = $target.Next()
$newPara1.Range.Text = "And this is also a proper way, with a different style"
$newPara1.Style = "FirstParagraph"

$r2 = $newPara1.Range
$r2.InsertParagraphAfter()
$newPara2 = $newPara1.Next()
$newPara2.Range.Text = "This is called the Intense Quote style."
$newPara2.Style = "BlockText"

Write-Host "done"
foreach ($p in $d.Paragraphs) {
    Write-Host ($p.Style.NameLocal + " | " + $p.Range.Text)
}
